# Auto-generated Excel COM-interop script
# Applies the "Added detailed docstring and comments" commit:
# adds new localisation rows (settings/export panel strings) to the
# Exiobase / Deutsch / English sheets, updates two existing German
# colormap translations, widens the English sheets column B, and
# leaves the "Deutsch" sheet as the active tab/selection.

$wb = $excel.ActiveWorkbook

$wsExiobase = $wb.Worksheets.Item("Exiobase")
$wsDeutsch  = $wb.Worksheets.Item("Deutsch")
$wsEnglish  = $wb.Worksheets.Item("English")

# ---- Exiobase sheet: append new key rows (A == B) ----
$wsExiobase.Range("A168").Value = "Open settings"
$wsExiobase.Range("B168").Value = "Open settings"
$wsExiobase.Rows.Item(168).RowHeight = 15
$wsExiobase.Range("A169").Value = "Impacts"
$wsExiobase.Range("B169").Value = "Impacts"
$wsExiobase.Rows.Item(169).RowHeight = 15
$wsExiobase.Range("A170").Value = "Include units in column names"
$wsExiobase.Range("B170").Value = "Include units in column names"
$wsExiobase.Rows.Item(170).RowHeight = 15
$wsExiobase.Range("A171").Value = "Localize column names"
$wsExiobase.Range("B171").Value = "Localize column names"
$wsExiobase.Rows.Item(171).RowHeight = 15
$wsExiobase.Range("A172").Value = "Output file"
$wsExiobase.Range("B172").Value = "Output file"
$wsExiobase.Rows.Item(172).RowHeight = 15
$wsExiobase.Range("A173").Value = "Excel Files (*.xlsx)"
$wsExiobase.Range("B173").Value = "Excel Files (*.xlsx)"
$wsExiobase.Rows.Item(173).RowHeight = 15
$wsExiobase.Range("A174").Value = "Please select at last one impact."
$wsExiobase.Range("B174").Value = "Please select at last one impact."
$wsExiobase.Rows.Item(174).RowHeight = 15
$wsExiobase.Range("A175").Value = "Please chose an .xlsx file."
$wsExiobase.Range("B175").Value = "Please chose an .xlsx file."
$wsExiobase.Rows.Item(175).RowHeight = 15
$wsExiobase.Range("A176").Value = "Excel export finished "
$wsExiobase.Range("B176").Value = "Excel export finished "
$wsExiobase.Rows.Item(176).RowHeight = 15
$wsExiobase.Range("A177").Value = "Failed to export Excel"
$wsExiobase.Range("B177").Value = "Failed to export Excel"
$wsExiobase.Rows.Item(177).RowHeight = 15

# ---- Deutsch sheet: fix two existing colormap translations ----
$wsDeutsch.Range("B146").Value = "Akzent "
$wsDeutsch.Range("B147").Value = "Dunkel 2"

# ---- Deutsch sheet: append new key/translation rows ----
$wsDeutsch.Range("A169").Value = "Open settings"
$wsDeutsch.Range("B169").Value = "Einstellungen öffnen"
$wsDeutsch.Rows.Item(169).RowHeight = 15
$wsDeutsch.Range("A170").Value = "Refresh"
$wsDeutsch.Range("B170").Value = "Aktualisieren"
$wsDeutsch.Rows.Item(170).RowHeight = 15
$wsDeutsch.Range("A171").Value = "Export data"
$wsDeutsch.Range("B171").Value = "Daten exportieren"
$wsDeutsch.Rows.Item(171).RowHeight = 15
$wsDeutsch.Range("A172").Value = "Impacts"
$wsDeutsch.Range("B172").Value = "Impacts"
$wsDeutsch.Rows.Item(172).RowHeight = 15
$wsDeutsch.Range("A173").Value = "Include units in column names"
$wsDeutsch.Range("B173").Value = "Einheiten in Spaltennamen aufnehmen"
$wsDeutsch.Rows.Item(173).RowHeight = 15
$wsDeutsch.Range("A174").Value = "Localize column names"
$wsDeutsch.Range("B174").Value = "Spaltennamen lokalisieren"
$wsDeutsch.Rows.Item(174).RowHeight = 15
$wsDeutsch.Range("A175").Value = "Output file"
$wsDeutsch.Range("B175").Value = "Ausgabedatei"
$wsDeutsch.Rows.Item(175).RowHeight = 15
$wsDeutsch.Range("A176").Value = "Excel Files (*.xlsx)"
$wsDeutsch.Range("B176").Value = "Excel-Dateien (*.xlsx)"
$wsDeutsch.Rows.Item(176).RowHeight = 15
$wsDeutsch.Range("A177").Value = "Please select at last one impact."
$wsDeutsch.Range("B177").Value = "Bitte mindestens einen Impact auswählen."
$wsDeutsch.Rows.Item(177).RowHeight = 15
$wsDeutsch.Range("A178").Value = "Please chose an .xlsx file."
$wsDeutsch.Range("B178").Value = "Bitte eine .xlsx-Datei wählen."
$wsDeutsch.Rows.Item(178).RowHeight = 15
$wsDeutsch.Range("A179").Value = "Excel export finished "
$wsDeutsch.Range("B179").Value = "Excel-Export abgeschlossen"
$wsDeutsch.Rows.Item(179).RowHeight = 15
$wsDeutsch.Range("A180").Value = "Failed to export Excel"
$wsDeutsch.Range("B180").Value = "Fehler beim Excel-Export"
$wsDeutsch.Rows.Item(180).RowHeight = 15

# ---- English sheet: append new key rows (A == B) ----
$wsEnglish.Range("A172").Value = "Open settings"
$wsEnglish.Range("B172").Value = "Open settings"
$wsEnglish.Rows.Item(172).RowHeight = 15
$wsEnglish.Range("A173").Value = "Impacts"
$wsEnglish.Range("B173").Value = "Impacts"
$wsEnglish.Rows.Item(173).RowHeight = 15
$wsEnglish.Range("A174").Value = "Include units in column names"
$wsEnglish.Range("B174").Value = "Include units in column names"
$wsEnglish.Rows.Item(174).RowHeight = 15
$wsEnglish.Range("A175").Value = "Localize column names"
$wsEnglish.Range("B175").Value = "Localize column names"
$wsEnglish.Rows.Item(175).RowHeight = 15
$wsEnglish.Range("A176").Value = "Output file"
$wsEnglish.Range("B176").Value = "Output file"
$wsEnglish.Rows.Item(176).RowHeight = 15
$wsEnglish.Range("A177").Value = "Excel Files (*.xlsx)"
$wsEnglish.Range("B177").Value = "Excel Files (*.xlsx)"
$wsEnglish.Rows.Item(177).RowHeight = 15
$wsEnglish.Range("A178").Value = "Please select at last one impact."
$wsEnglish.Range("B178").Value = "Please select at last one impact."
$wsEnglish.Rows.Item(178).RowHeight = 15
$wsEnglish.Range("A179").Value = "Please chose an .xlsx file."
$wsEnglish.Range("B179").Value = "Please chose an .xlsx file."
$wsEnglish.Rows.Item(179).RowHeight = 15
$wsEnglish.Range("A180").Value = "Excel export finished "
$wsEnglish.Range("B180").Value = "Excel export finished "
$wsEnglish.Rows.Item(180).RowHeight = 15
$wsEnglish.Range("A181").Value = "Failed to export Excel"
$wsEnglish.Range("B181").Value = "Failed to export Excel"
$wsEnglish.Rows.Item(181).RowHeight = 15

# ---- English sheet: widen column B to fit the longer strings ----
$wsEnglish.Columns.Item(2).ColumnWidth = 114.16666666666667

# ---- Restore per-sheet selections ----
$wsExiobase.Range("B179").Select() | Out-Null
$wsEnglish.Range("B174:B181").Select() | Out-Null

# ---- Deutsch becomes the active/selected tab ----
$wsDeutsch.Activate() | Out-Null
$wsDeutsch.Range("A173:A180").Select() | Out-Null

Write-Host "edit.ps1 completed"
